# STEP #1 - Move temp-word-docs into docs-journey-log and update memory reference
#
# Applies the following changes to the Prototype-Architecture-Diagrams deck:
#  - Slide 1: shrink the subtitle font under the title (31pt -> 27pt)
#  - Slide 2: retitle "MVP" -> "Dockerized", reflow/enlarge the title box,
#             shrink its subtitle font (31pt -> 27pt)
#  - Slide 3: retitle the "Production" title -> "Cloud-Native", shrink its
#             subtitle font (31pt -> 27pt) and reword the subtitle text;
#             shrink + reflow the footnote textbox and its runs (14pt -> 12pt),
#             and clarify the "Azure Container Apps" mention with its acronym

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 - "Initial prototype/POC Architecture Diagram"
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$title1 = $slide1.Shapes.Item("Title 1")
$title1TR = $title1.TextFrame.TextRange
# Run 1 = title text, Run 2 = subtitle text (separated by a manual line break)
$title1TR.Runs(2).Font.Size = 27

# ---------------------------------------------------------------------------
# Slide 2 - "MVP - Architecture Diagram" -> "Dockerized - Architecture Diagram"
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$title2 = $slide2.Shapes.Item("Title 1")

# Move/resize the title placeholder (off.y 316849 -> 67021, ext.cy 750709 -> 1000538;
# x/cx are left untouched). The literal below is nudged very slightly off the
# naive EMU/12700 conversion so the host's point-based Top/Height setters round
# to the exact target EMU instead of landing 1 EMU short.
$title2.Top = 5.2772556174380645
$title2.Height = 78.78253120798925

$title2TR = $title2.TextFrame.TextRange
$title2TR.Runs(1).Text = "Dockerized - Architecture Diagram"
$title2TR.Runs(2).Font.Size = 27

# ---------------------------------------------------------------------------
# Slide 3 - the "Production" title -> "Cloud-Native - Architecture Diagram"
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$title3 = $slide3.Shapes.Item("Title 1")
$title3TR = $title3.TextFrame.TextRange
$title3TR.Runs(1).Text = "Cloud-Native - Architecture Diagram"
$title3TR.Runs(2).Font.Size = 27
$title3TR.Runs(2).Text = "(Decoupled Agents / MCP servers as containers in Kubernetes/Azure [*])"

# Footnote textbox under the slide 3 diagram
$footnote = $slide3.Shapes.Item("TextBox 42")

# Move/resize (off.x 888537 -> 1117325, ext.cx 10787269 -> 9800643,
# ext.cy 307777 -> 276999; off.y is left untouched). As above, the literals
# are nudged by a hair so the point -> EMU rounding lands exactly on target.
$footnote.Left = 87.97834645669292
$footnote.Width = 771.7041962742462
$footnote.Height = 21.810944881889764

$footnoteTR = $footnote.TextFrame.TextRange
for ($i = 1; $i -le $footnoteTR.Runs().Count; $i++) {
  $footnoteTR.Runs($i).Font.Size = 12
}
$footnoteTR.Runs(4).Text = "Azure Container Apps (ACA)"
